$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (changed) date column (C) for data rows 2-14
# from serial date 45177 (2023-09-08) to 45178 (2023-09-09).
for ($row = 2; $row -le 14; $row++) {
    $ws.Cells.Item($row, 3).Value = 45178
}
